$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: "Cake Slice" -> "Cake Slice Varieties", drop the Ingredients value,
# and replace the Allergens text with the new wording. Row grows taller to
# fit the new wrapped allergens text (matches the existing styling used by
# other wrapped rows, e.g. row 6).
$ws.Range("B8").ClearContents() | Out-Null
$ws.Range("C8").Value = "Wheat, milk, eggs, gluten, soy. May contain peanuts, sesame."
$ws.Range("A8").Value = "Cake Slice Varieties"
$ws.Rows.Item(8).RowHeight = 30

# Update the current selection to match where editing left off.
$ws.Range("A8").Select() | Out-Null
